# Add 2022-Q3 data:
#  - insert a brand-new "2022-Q3" sheet right after "总计", pushing the
#    existing quarter sheets (2022-Q2, 2022-Q1, 2021-Q3, 2021-Q2) one
#    position later (their own name/data/formatting travel with them,
#    unchanged);
#  - populate "2022-Q3" with the new per-fund figures (copying the
#    "2022-Q2" sheet first so header/row formatting matches the other
#    quarter tabs);
#  - extend the "总计" (totals) sheet with the corresponding new row,
#    shifting the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# --- build the new "2022-Q3" sheet from a copy of "2022-Q2" (keeps the
#     same header row / column styling as every other quarter tab) ---
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# code / name / size / total position / position share / held value / rank
$wsQ3.Cells.Item(2,2).Formula = "'118002"
$wsQ3.Cells.Item(2,3).Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$wsQ3.Cells.Item(2,4).Formula = "'1.85"
$wsQ3.Cells.Item(2,5).Formula = "'93.04"
$wsQ3.Cells.Item(2,6).Formula = "'9.65"
$wsQ3.Cells.Item(2,7).Formula = "'0.1785"
$wsQ3.Cells.Item(2,8).Value = 2

$wsQ3.Cells.Item(3,2).Formula = "'000593"
$wsQ3.Cells.Item(3,3).Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$wsQ3.Cells.Item(3,4).Formula = "'1.85"
$wsQ3.Cells.Item(3,5).Formula = "'93.04"
$wsQ3.Cells.Item(3,6).Formula = "'9.65"
$wsQ3.Cells.Item(3,7).Formula = "'0.1785"
$wsQ3.Cells.Item(3,8).Value = 2

$wsQ3.Cells.Item(4,2).Formula = "'005676"
$wsQ3.Cells.Item(4,3).Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$wsQ3.Cells.Item(4,4).Formula = "'1.85"
$wsQ3.Cells.Item(4,5).Formula = "'93.04"
$wsQ3.Cells.Item(4,6).Formula = "'9.65"
$wsQ3.Cells.Item(4,7).Formula = "'0.1785"
$wsQ3.Cells.Item(4,8).Value = 2

# --- extend the "总计" sheet: insert the 2022-Q3 row at the top of the
#     data and push the rest down by one row ---
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 3
$wsTotal.Cells.Item(2,4).Value = 0.54

$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(4,1).Value = 2
$wsTotal.Cells.Item(5,1).Value = 3
$wsTotal.Cells.Item(6,1).Value = 4
